$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = "'43.910.04"
$ws.Range("E2").Value = "  +5.26%  "
$ws.Range("D3").Value2 = "'2.294.52"
$ws.Range("E3").Value = "  +2.97%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value2 = "'231.47"
$ws.Range("E5").Value = "  +0.20%  "
$ws.Range("E6").Value = "  +0.98%  "
$ws.Range("D7").Value2 = "'63.75"
$ws.Range("E7").Value = "  +5.68%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("E9").Value = "  +5.00%  "
$ws.Range("D10").Value2 = "'0.0960"
$ws.Range("E10").Value = "  +6.90%  "
$ws.Range("D11").Value2 = "'57.65"
$ws.Range("E11").Value = "  -0.77%  "
$ws.Range("D12").Value2 = "'26.40"
$ws.Range("E12").Value = "  +14.99%  "
$ws.Range("E13").Value = "  +0.54%  "
$ws.Range("D14").Value2 = "'2.635.82"
$ws.Range("E14").Value = "  +2.79%  "
$ws.Range("D15").Value2 = "'15.94"
$ws.Range("E15").Value = "  +1.80%  "
$ws.Range("E16").Value = "  +6.66%  "
$ws.Range("E17").Value = "  +1.97%  "
$ws.Range("D18").Value2 = "'2.292.05"
$ws.Range("E18").Value = "  +2.28%  "
$ws.Range("D19").Value2 = "'43.821.61"
$ws.Range("E19").Value = "  +5.11%  "
$ws.Range("D20").Value2 = "'0.0₃0982"
$ws.Range("E20").Value = "  +9.22%  "
$ws.Range("D21").Value2 = "'73.59"
$ws.Range("E21").Value = "  +2.00%  "
$ws.Range("D22").Value2 = "'6.23"
$ws.Range("E22").Value = "  +2.11%  "
$ws.Range("D23").Value2 = "'254.66"
$ws.Range("E23").Value = "  +2.87%  "
$ws.Range("D24").Value2 = "'2.65"
$ws.Range("E24").Value = "  +12.26%  "
$ws.Range("E25").Value = "  +0.00%  "
$ws.Range("D26").Value2 = "'2.32"
$ws.Range("E26").Value = "  -2.44%  "
$ws.Range("D27").Value2 = "'9.92"
$ws.Range("E27").Value = "  +2.34%  "
$ws.Range("D28").Value2 = "'171.60"
$ws.Range("E28").Value = "  +1.34%  "
$ws.Range("E29").Value = "  -1.15%  "
$ws.Range("D30").Value2 = "'20.68"
$ws.Range("E30").Value = "  +4.04%  "
$ws.Range("E31").Value = "  +2.65%  "
$ws.Range("D32").Value2 = "'2.77"
$ws.Range("E32").Value = "  +4.20%  "
$ws.Range("D33").Value2 = "'0.123"
$ws.Range("E33").Value = "  +0.82%  "
$ws.Range("D34").Value2 = "'0.0690"
$ws.Range("E34").Value = "  +6.68%  "
$ws.Range("D35").Value2 = "'5.12"
$ws.Range("E35").Value = "  +2.37%  "
$ws.Range("D36").Value2 = "'4.76"
$ws.Range("E36").Value = "  +2.20%  "
$ws.Range("D37").Value2 = "'3.72"
$ws.Range("E37").Value = "  +2.57%  "
$ws.Range("E38").Value = "  +0.36%  "
$ws.Range("E39").Value = "  +0.44%  "
$ws.Range("E40").Value = "  +4.55%  "
$ws.Range("E41").Value = "  +0.12%  "
$ws.Range("D42").Value2 = "'8.79"
$ws.Range("E42").Value = "  +2.56%  "
$ws.Range("D43").Value2 = "'11.02"
$ws.Range("E43").Value = "  +29.69%  "
$ws.Range("D44").Value2 = "'0.000221"
$ws.Range("E44").Value = "  -4.87%  "
$ws.Range("D45").Value2 = "'4.56"
$ws.Range("E45").Value = "  +0.60%  "
$ws.Range("D46").Value2 = "'1.23"
$ws.Range("E46").Value = "  +1.07%  "
$ws.Range("D47").Value2 = "'99.18"
$ws.Range("E47").Value = "  +0.68%  "
$ws.Range("D48").Value2 = "'0.0969"
$ws.Range("E48").Value = "  +0.95%  "
$ws.Range("D49").Value2 = "'17.29"
$ws.Range("E49").Value = "  +4.35%  "
$ws.Range("D50").Value2 = "'1.488.94"
$ws.Range("E50").Value = "  +1.83%  "
$ws.Range("D51").Value2 = "'2.32"
$ws.Range("E51").Value = "  +0.88%  "
